$d = $word.ActiveDocument

function Replace-InParagraph($index, $old, $new) {
    $p = $d.Paragraphs($index)
    $r = $p.Range
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Replace failed in paragraph $index : '$old' -> '$new'"
    }
}

# --- Summary stats block (elapsed time / speed) ---
Replace-InParagraph 6 "Elapsed Time: 17.66s" "Elapsed Time: 18.83s"
Replace-InParagraph 7 "Speed: 0.17 pages/sec" "Speed: 0.16 pages/sec"

# --- Category section (was "Books"/"A Light in the ...") now "Horror"/"Security" ---
Replace-InParagraph 18 "Books" "Horror"
Replace-InParagraph 19 "https://books.toscrape.com/catalogue/category/books_1/index.html" "https://books.toscrape.com/catalogue/category/books/horror_31/index.html"
Replace-InParagraph 20 "Home" "Home > Books"
Replace-InParagraph 22 "A Light in the ..." "Security"
Replace-InParagraph 24 "A Light in the ..." "Security"
Replace-InParagraph 24 "£51.77" "£39.25"

# --- Category section (was "Fantasy"/"Unicorn Tracks") now "Thriller"/"In Her Wake" ---
Replace-InParagraph 26 "Fantasy" "Thriller"
Replace-InParagraph 27 "https://books.toscrape.com/catalogue/category/books/fantasy_19/index.html" "https://books.toscrape.com/catalogue/category/books/thriller_37/index.html"
Replace-InParagraph 30 "Unicorn Tracks" "In Her Wake"
Replace-InParagraph 32 "Unicorn Tracks" "In Her Wake"
Replace-InParagraph 32 "£18.78" "£12.84"
